$d = $word.ActiveDocument
$r = $d.Content
$found = $r.Find.Execute("PT PLN (PERSERO) ", $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter("KANTOR ")
$r.Select()
$word.Selection.LanguageID = 1057
Write-Output ("Sel Lang=" + $word.Selection.LanguageID)
Write-Output ("Range Lang=" + $r.LanguageID)
